$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that currently sits at the end of the
#    "View your account details." paragraph. It gets re-created later at the
#    new final location (end of the "Role: Driver" paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Website task: "Log into the website using these details:" used to be
#    followed by a single empty paragraph. Turn it into two paragraphs
#    carrying the login credentials.
# ---------------------------------------------------------------------------
$idPara = $d.Paragraphs.Item(19)
$idPara.Range.Text = "ID: A3456"
$idPara.Range.InsertParagraphAfter()

$pwPara = $d.Paragraphs.Item(20)
$pwPara.Range.Text = "Password: coachesrule1"

# ---------------------------------------------------------------------------
# 3) Website task: "Create a new employee with these details:" used to be
#    followed by an empty (List Paragraph styled) paragraph and then a plain
#    empty paragraph. Replace those two with five indented paragraphs
#    describing the new employee, moving the "_GoBack" bookmark onto the end
#    of the last one ("Role: Driver").
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(24)
$p1.Style = "Normal"
$p1.LeftIndent = 18
$p1.Range.Text = "ID: D4653"
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(25)
$p2.Style = "Normal"
$p2.LeftIndent = 18
$p2.Range.Text = "Password: driver123"
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(26)
$p3.Style = "Normal"
$p3.LeftIndent = 18
$p3.Range.Text = "First Name: Joe"
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(27)
$p4.Style = "Normal"
$p4.LeftIndent = 18
$p4.Range.Text = "Last Name: Smith"
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs.Item(28)
$p5.Style = "Normal"
$p5.LeftIndent = 18
$p5.Range.Text = "Role: Driver"

$bmRange = $d.Range($p5.Range.End - 1, $p5.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
